$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
# Insert new "Priority (0-100)" / "Weight (0-100)" numeric-range headers in
# B1/C1, pushing the old Project/Earliest Date/Due Date/Notes headers one
# column to the right (D..G).
$ws.Range("G1").Value = "Notes"
$ws.Range("F1").Value = "Due Date"
$ws.Range("E1").Value = "Earliest Date"
$ws.Range("D1").Value = "Project"
$ws.Range("B1").Value = "Priority (0-100)"
$ws.Range("C1").Value = "Weight (0-100)"

# --- Rows 2 & 3 ---------------------------------------------------------
# Only the Priority/Weight numbers change (rescaled from a 1-10 scale to a
# 0-100 scale); Name/Project/Dates/Notes stay as-is.
$ws.Range("B2").Value = 40
$ws.Range("C2").Value = 70

$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 100

# --- Column widths ----------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 21.5703125
$ws.Columns.Item(3).ColumnWidth = 19.7109375

# --- Selected cell ------------------------------------------------------------
$ws.Range("C7").Select() | Out-Null
